{"js": "// Update the two model-summary tables (\"chisq\" ANOVA tables) with refreshed\n// values from the model run that adds random effects.\n\n// Helper: replace a table cell's text in place (keeps the existing run /\n// paragraph formatting, just swaps the literal text) instead of clobbering\n// the cell body, which would lose formatting like xml:space=\"preserve\".\nfunction setCellText(cell, text) {\n  const range = cell.body.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// ---- Table 1 (richness model) ----\nconst t1 = tables.items[0];\nt1.rows.load(\"items\");\nawait context.sync();\nt1.rows.items.forEach((r) => r.cells.load(\"items\"));\nawait context.sync();\n\n// Row 1: (Intercept) -> chisq 58.49 -> 68.36\nsetCellText(t1.rows.items[1].cells.items[1], \"68.36\");\n// Row 2: poly(elevation_mean, 2) -> chisq 78.63 -> 168.90\nsetCellText(t1.rows.items[2].cells.items[1], \"168.90\");\n// Row 3: regions -> chisq 15.98 -> 21.66\nsetCellText(t1.rows.items[3].cells.items[1], \"21.66\");\n// Row 4: seasons -> chisq 2.13 -> 0.21, pr_chisq 0.14 -> 0.65\nsetCellText(t1.rows.items[4].cells.items[1], \"0.21\");\nsetCellText(t1.rows.items[4].cells.items[3], \"0.65\");\nawait context.sync();\n\n// Remove rows \"poly(elevation_mean, 2):seasons\" and \"regions:seasons\" (now\n// folded into the random-effects structure instead of fixed interactions).\nt1.rows.load(\"items\");\nawait context.sync();\nt1.rows.items[5].cells.load(\"items\");\nt1.rows.items[6].cells.load(\"items\");\nawait context.sync();\n// Delete from the bottom up so earlier indices stay valid.\nt1.rows.items[6].delete();\nt1.rows.items[5].delete();\nawait context.sync();\n\n// ---- Table 2 (occurrence model) ----\nconst t2 = tables.items[1];\nt2.rows.load(\"items\");\nawait context.sync();\nt2.rows.items.forEach((r) => r.cells.load(\"items\"));\nawait context.sync();\n\n// Row 1: (Intercept) -> chisq 43.54 -> 27.70\nsetCellText(t2.rows.items[1].cells.items[1], \"27.70\");\n// Row 2: poly(elevation_mean, 2) -> chisq 70.34 -> 56.48\nsetCellText(t2.rows.items[2].cells.items[1], \"56.48\");\n// Row 3: regions -> chisq 1.96 -> 1.55, pr_chisq 0.38 -> 0.46\nsetCellText(t2.rows.items[3].cells.items[1], \"1.55\");\nsetCellText(t2.rows.items[3].cells.items[3], \"0.46\");\n// Row 4: seasons -> chisq 0.50 -> 0.46, pr_chisq 0.48 -> 0.50\nsetCellText(t2.rows.items[4].cells.items[1], \"0.46\");\nsetCellText(t2.rows.items[4].cells.items[3], \"0.50\");\n// Row 5: poly(elevation_mean, 2):regions -> chisq 6.94 -> 4.88, pr_chisq 0.14 -> 0.30\nsetCellText(t2.rows.items[5].cells.items[1], \"4.88\");\nsetCellText(t2.rows.items[5].cells.items[3], \"0.30\");\n// Row 6: poly(elevation_mean, 2):seasons -> chisq 7.73 -> 10.84, pr_chisq 0.02 -> 0.00, signif * -> **\nsetCellText(t2.rows.items[6].cells.items[1], \"10.84\");\nsetCellText(t2.rows.items[6].cells.items[3], \"0.00\");\nsetCellText(t2.rows.items[6].cells.items[4], \"**\");\n// Row 7: regions:seasons -> chisq 3.81 -> 2.42, pr_chisq 0.15 -> 0.30\nsetCellText(t2.rows.items[7].cells.items[1], \"2.42\");\nsetCellText(t2.rows.items[7].cells.items[3], \"0.30\");\n// Row 8: poly(elevation_mean, 2):regions:seasons -> chisq 23.67 -> 24.57\nsetCellText(t2.rows.items[8].cells.items[1], \"24.57\");\n\nawait context.sync();\n", "ps1": "# Update the two model-summary tables (\"chisq\" ANOVA tables) with refreshed\n# values from the model run that adds random effects.\n$d = $word.ActiveDocument\n\n# ---- Table 1 (richness model) ----\n$t1 = $d.Tables.Item(1)\n\n# Row 2: (Intercept) -> chisq 58.49 -> 68.36\n$t1.Cell(2,2).Range.Text = \"68.36\"\n# Row 3: poly(elevation_mean, 2) -> chisq 78.63 -> 168.90\n$t1.Cell(3,2).Range.Text = \"168.90\"\n# Row 4: regions -> chisq 15.98 -> 21.66\n$t1.Cell(4,2).Range.Text = \"21.66\"\n# Row 5: seasons -> chisq 2.13 -> 0.21, pr_chisq 0.14 -> 0.65\n$t1.Cell(5,2).Range.Text = \"0.21\"\n$t1.Cell(5,4).Range.Text = \"0.65\"\n\n# Remove rows \"poly(elevation_mean, 2):seasons\" (row 6) and \"regions:seasons\"\n# (row 7) - folded into the random-effects structure instead of fixed\n# interactions. Delete from the bottom up so earlier row indices stay valid.\n$t1.Rows.Item(7).Delete()\n$t1.Rows.Item(6).Delete()\n\n# ---- Table 2 (occurrence model) ----\n$t2 = $d.Tables.Item(2)\n\n# Row 2: (Intercept) -> chisq 43.54 -> 27.70\n$t2.Cell(2,2).Range.Text = \"27.70\"\n# Row 3: poly(elevation_mean, 2) -> chisq 70.34 -> 56.48\n$t2.Cell(3,2).Range.Text = \"56.48\"\n# Row 4: regions -> chisq 1.96 -> 1.55, pr_chisq 0.38 -> 0.46\n$t2.Cell(4,2).Range.Text = \"1.55\"\n$t2.Cell(4,4).Range.Text = \"0.46\"\n# Row 5: seasons -> chisq 0.50 -> 0.46, pr_chisq 0.48 -> 0.50\n$t2.Cell(5,2).Range.Text = \"0.46\"\n$t2.Cell(5,4).Range.Text = \"0.50\"\n# Row 6: poly(elevation_mean, 2):regions -> chisq 6.94 -> 4.88, pr_chisq 0.14 -> 0.30\n$t2.Cell(6,2).Range.Text = \"4.88\"\n$t2.Cell(6,4).Range.Text = \"0.30\"\n# Row 7: poly(elevation_mean, 2):seasons -> chisq 7.73 -> 10.84, pr_chisq 0.02 -> 0.00, signif * -> **\n$t2.Cell(7,2).Range.Text = \"10.84\"\n$t2.Cell(7,4).Range.Text = \"0.00\"\n$t2.Cell(7,5).Range.Text = \"**\"\n# Row 8: regions:seasons -> chisq 3.81 -> 2.42, pr_chisq 0.15 -> 0.30\n$t2.Cell(8,2).Range.Text = \"2.42\"\n$t2.Cell(8,4).Range.Text = \"0.30\"\n# Row 9: poly(elevation_mean, 2):regions:seasons -> chisq 23.67 -> 24.57\n$t2.Cell(9,2).Range.Text = \"24.57\"\n"}
